# Update "Pais" sheet: refresh the timestamp caption and a batch of country
# statistics, plus re-order three pairs of countries whose updated "Casos
# totales" (column B) changed their rank in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row {
    param($row, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Update the "last refreshed" caption in A1 -----------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 13:26"

# --- Straightforward numeric refreshes (country stays on the same row) ----
# Iran (row 15)
Set-Row 15 402029 2089 346242 32630 0 128 23157

# Oman (row 40)
Set-Row 40 89746 1409 83771 5195 0 18 780

# Emiratos Arabes Unidos (row 46)
Set-Row 46 79489 640 69451 9639 0 0 399

# Estado de Palestina (row 73)
Set-Row 73 30574 668 20082 10271 0 11 221

# Libia (row 79)
Set-Row 79 22781 433 12183 10236 0 8 362

# Madagascar (row 85)
Set-Row 85 15757 20 14368 1178 0 1 211

# Senegal (row 87)
Set-Row 87 14280 43 10520 3463 0 2 297

# Hong Kong (row 116)
Set-Row 116 4958 19 4630 228 0 0 100

# Islandia (row 150)
Set-Row 150 2165 3 2092 63 0 0 10

# Islas Feroe (row 179)
Set-Row 179 423 5 410 13 0 0 0

# Gibraltar (row 183)
Set-Row 183 330 3 298 32 0 0 0

# --- Re-ranked pairs: updated country now outranks its former neighbour ---
# Rumania jumps ahead of Republica Dominicana (rows 34/35)
$ws.Range("A34").Value = "Rumania"
Set-Row 34 103495 1109 43025 56307 0 36 4163
$ws.Range("A35").Value = "Republica Dominicana"
Set-Row 35 103092 0 76531 24608 0 0 1953

# Suiza jumps ahead of Uzbekistan (rows 60/61)
$ws.Range("A60").Value = "Suiza"
Set-Row 60 47179 475 38900 6258 0 1 2021
$ws.Range("A61").Value = "Uzbekistan"
Set-Row 61 46850 129 43511 2953 0 2 386

# Malta jumps ahead of Guinea-Bisau (rows 145/146)
$ws.Range("A145").Value = "Malta"
Set-Row 145 2352 78 1872 465 0 0 15
$ws.Range("A146").Value = "Guinea-Bisau"
Set-Row 146 2275 0 1127 1109 0 0 39
